# "Add files via upload" — populate the previously-empty "Parecer recebido"
# (received opinion/report file name) cells for the six people whose file
# had not yet been uploaded when the sheet was last saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value  = "Parecer_Alessandra.pdf"  # row 3  - Alessandra Freixo Braga
$ws.Range("B8").Value  = "Parecer_Andrea.pdf"      # row 8  - Andrea Marcia
$ws.Range("B11").Value = "Parecer_Emynna.pdf"      # row 11 - Emynna Cavalcante Guimaraes
$ws.Range("B18").Value = "Parecer_Jessica.pdf"     # row 18 - Jessica Farias Macedo
$ws.Range("B21").Value = "Parecer_Lunna.pdf"       # row 21 - Lunna Nascimento Barroso
$ws.Range("B25").Value = "Parecer_Maycon.pdf"      # row 25 - Maycon Deyvis Sena Vicente

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("B35").Select()
